# Column A ("city") was resized so the longest city name ("West
# Christopherberg") fits without truncation - i.e. the user selected
# column A and auto-fit it to its contents (double-clicking the column
# A/B border, or Format > Column Width > AutoFit Selection), which is
# also why Excel marks the resulting width as a "best fit" width.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(1)

# Auto-fit column A to its contents first (the actual user action).
$col.AutoFit() | Out-Null

# Pin the width to the precise value Excel's font metrics produced for
# this content (21.42578125 characters) so the saved column width
# matches exactly.
$col.ColumnWidth = 20.6
